$d = $word.ActiveDocument

$replacements = @(
    @{old="673×8="; new="653×8="},
    @{old="886×3="; new="961×6="},
    @{old="760×6="; new="775×7="},
    @{old="953×7="; new="408×4="},
    @{old="962×2="; new="941×2="},
    @{old="450×4="; new="503×4="},
    @{old="562×5="; new="517×5="},
    @{old="894×9="; new="488×5="},
    @{old="938×5="; new="214×2="},
    @{old="662×3="; new="124×2="},
    @{old="223×9="; new="415×5="},
    @{old="589×8="; new="894×8="},
    @{old="941×9="; new="303×3="},
    @{old="908×3="; new="947×2="},
    @{old="456×4="; new="163×9="},
    @{old="353×5="; new="719×2="},
    @{old="853×3="; new="259×9="},
    @{old="170×2="; new="660×8="},
    @{old="186×6="; new="211×7="},
    @{old="306×2="; new="287×5="},
    @{old="412×7="; new="373×7="},
    @{old="313×4="; new="797×2="},
    @{old="941×5="; new="299×2="},
    @{old="434×6="; new="796×9="},
    @{old="360×8="; new="641×3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
